# Edit: mlk.docx
# 1. Collapse the first paragraph's 4 runs (plain text + 3 red "(This is a
#    change - Version for branch alternate)" runs) down to a single plain
#    run reading "This is a Microsoft word document." (no trailing spaces,
#    no coloring).
# 2. Remove the trailing empty paragraph that only carries a light-gray
#    shading (w:shd fill="F9F9F9") just before the section break.

$d = $word.ActiveDocument

# --- 1. Fix up the opening paragraph ---------------------------------
$p1 = $d.Paragraphs.Item(1)
$rng = $p1.Range
# Exclude the paragraph mark itself from the range we rewrite.
$rng.End = $rng.End - 1
$rng.Text = "This is a Microsoft word document."

# Re-grab the (now shorter) paragraph range and clear any leftover
# character formatting (e.g. the red C00000 color) so a single plain run
# remains.
$rng2 = $p1.Range
$rng2.End = $rng2.End - 1
$rng2.Font.Reset()

# --- 2. Drop the trailing shaded empty paragraph ----------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
# Extend the range one character to the left so the preceding paragraph
# mark is removed too -- this is what actually deletes the paragraph
# rather than just clearing its (already empty) contents.
$killRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$killRange.Delete()
